# Updates cryptocurrency price/volume data to reflect the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.549.63"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").Value = "2.295.54"
$ws.Range("E3").Value = "  +0.52%  "
$ws.Range("E4").Value = "  +0.72%  "
$ws.Range("D5").Value = "'312.39"
$ws.Range("E5").Value = "  -1.95%  "
$ws.Range("D6").Value = "104.42"
$ws.Range("E6").Value = "  +3.48%  "
$ws.Range("E7").Value = "  -0.42%  "
$ws.Range("E8").Value = "  -0.23%  "
$ws.Range("D9").Value = "'0.602"
$ws.Range("E9").Value = "  +0.12%  "
$ws.Range("D10").Value = "'39.11"
$ws.Range("E10").Value = "  +0.12%  "
$ws.Range("D11").Value = "'0.0903"
$ws.Range("E11").Value = "  +0.28%  "
$ws.Range("D12").Value = "'8.26"
$ws.Range("E12").Value = "  +0.39%  "
$ws.Range("E13").Value = "  +1.40%  "
$ws.Range("D14").Value = "'0.987"
$ws.Range("E14").Value = "  +3.29%  "
$ws.Range("D15").Value = "15.15"
$ws.Range("D16").Value = "2.644.73"
$ws.Range("E16").Value = "  +0.38%  "
$ws.Range("D17").Value = "2.290.94"
$ws.Range("E17").Value = "  +0.10%  "
$ws.Range("D18").Value = "42.731.25"
$ws.Range("E18").Value = "  +1.11%  "
$ws.Range("D19").Value = "'7.30"
$ws.Range("E19").Value = "  -0.91%  "
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("D21").Value = "'13.53"
$ws.Range("E21").Value = "  +4.40%  "
$ws.Range("D22").Value = "'73.33"
$ws.Range("E22").Value = "  +0.90%  "
$ws.Range("D23").Value = "3.48"
$ws.Range("E23").Value = "  -1.95%  "
$ws.Range("D24").Value = "'263.96"
$ws.Range("E24").Value = "  -0.99%  "
$ws.Range("D25").Value = "'2.18"
$ws.Range("E25").Value = "  -0.93%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  -0.16%  "
$ws.Range("D27").Value = "10.74"
$ws.Range("E27").Value = "  -0.39%  "
$ws.Range("D28").Value = "'7.09"
$ws.Range("E28").Value = "  +16.45%  "
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("D30").Value = "'22.39"
$ws.Range("E30").Value = "  -0.14%  "
$ws.Range("D31").Value = "35.93"
$ws.Range("E31").Value = "  -3.59%  "
$ws.Range("D32").Value = "'165.15"
$ws.Range("E32").Value = "  -0.29%  "
$ws.Range("D33").Value = "'0.0863"
$ws.Range("E33").Value = "  -0.71%  "
$ws.Range("E34").Value = "  -1.82%  "
$ws.Range("E35").Value = "  +0.29%  "
$ws.Range("D36").Value = "0.113"
$ws.Range("E36").Value = "  -1.55%  "
$ws.Range("D37").Value = "'4.50"
$ws.Range("E37").Value = "  -1.37%  "
$ws.Range("E38").Value = "  -1.56%  "
$ws.Range("E39").Value = "  +2.87%  "
$ws.Range("D40").Value = "'2.74"
$ws.Range("E40").Value = "  +0.32%  "
$ws.Range("D41").Value = "'1.60"
$ws.Range("E41").Value = "  +5.34%  "
$ws.Range("D42").Value = "'99.37"
$ws.Range("E42").Value = "  +7.06%  "
$ws.Range("D43").Value = "'69.26"
$ws.Range("E43").Value = "  +1.27%  "
$ws.Range("D44").Value = "'0.227"
$ws.Range("E44").Value = "  +1.67%  "
$ws.Range("D45").Value = "1.01"
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("D48").Value = "'79.50"
$ws.Range("E48").Value = "  +1.15%  "
$ws.Range("D49").Value = "'110.99"
$ws.Range("E49").Value = "  -3.02%  "
$ws.Range("D50").Value = "'5.19"
$ws.Range("E50").Value = "  -0.21%  "
$ws.Range("D51").Value = "'8.67"
$ws.Range("E51").Value = "  -2.89%  "

# Rows 46 and 47 swapped coins: Celestia and Maker traded ranking positions.
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "1.750.25"
$ws.Range("E46").Value = "  +9.58%  "
$ws.Range("B47").Value = "Celestia"
$ws.Range("C47").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D47").Value = "'12.10"
$ws.Range("E47").Value = "  +1.85%  "
